$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column "29-nov" before the
#     existing "01-oct." column (at EB), shifting everything after it
#     one column to the right (through FG).
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Columns("EB:EB").Insert()
$wsPrix.Range("EB1").Value = "29-nov"
for ($r = 2; $r -le 25; $r++) {
    $wsPrix.Cells.Item($r, 132).Value = "-"
}

# --- Sheet "Gaz": append the next day's price row.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A162").Value = "'2025-11-27"
$wsGaz.Range("B162").Value = 27.875

# --- Sheet "CO2": append the next day's price row.
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A162").Value = "'2025-11-27"
$wsCo2.Range("B162").Value = 82.3
